# Update cryptos list (price/volume) values, preserving text cell type/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" "65.242.39"
Set-TextCell "E2" "  -1.98%  "
Set-TextCell "D3" "3.483.43"
Set-TextCell "E3" "  -1.13%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  +0.00%  "
Set-TextCell "D5" "588.88"
Set-TextCell "E5" "  -2.92%  "
Set-TextCell "D6" "137.40"
Set-TextCell "E6" "  -4.47%  "
Set-TextCell "D7" "3.482.13"
Set-TextCell "E7" "  -1.15%  "
Set-TextCell "E8" "  +0.10%  "
Set-TextCell "D9" "0.491"
Set-TextCell "E9" "  -3.09%  "
Set-TextCell "D10" "0.123"
Set-TextCell "E10" "  -5.81%  "
Set-TextCell "D11" "7.18"
Set-TextCell "E11" "  -6.98%  "
Set-TextCell "D12" "0.383"
Set-TextCell "E12" "  -5.63%  "
Set-TextCell "D13" "4.069.41"
Set-TextCell "E13" "  -1.06%  "
Set-TextCell "D14" "0.0000183"
Set-TextCell "E14" "  -6.34%  "
Set-TextCell "B15" "Avalanche"
Set-TextCell "C15" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D15" "26.56"
Set-TextCell "E15" "  -7.45%  "
Set-TextCell "B16" "WrappedEther"
Set-TextCell "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D16" "3.467.96"
Set-TextCell "E16" "  -2.04%  "
Set-TextCell "E17" "  -1.21%  "
Set-TextCell "D18" "65.125.88"
Set-TextCell "E18" "  -1.87%  "
Set-TextCell "D19" "9.73"
Set-TextCell "E19" "  -9.61%  "
Set-TextCell "D20" "5.78"
Set-TextCell "E20" "  -5.76%  "
Set-TextCell "D21" "13.93"
Set-TextCell "E21" "  -4.60%  "
Set-TextCell "D22" "389.40"
Set-TextCell "E22" "  -7.91%  "
Set-TextCell "D23" "0.555"
Set-TextCell "E23" "  -5.83%  "
Set-TextCell "E24" "  +0.05%  "
Set-TextCell "D25" "72.55"
Set-TextCell "E25" "  -5.82%  "
Set-TextCell "E26" "  -0.07%  "
Set-TextCell "D27" "3.621.82"
Set-TextCell "E27" "  -1.24%  "
Set-TextCell "E28" "  -4.05%  "
Set-TextCell "D29" "0.998"
Set-TextCell "E29" "  -0.23%  "
Set-TextCell "E30" "  -6.62%  "
Set-TextCell "D31" "8.23"
Set-TextCell "E31" "  -7.68%  "
Set-TextCell "D32" "2.22"
Set-TextCell "E32" "  -9.84%  "
Set-TextCell "D33" "3.500.01"
Set-TextCell "E33" "  -0.79%  "
Set-TextCell "E34" "  -0.02%  "
Set-TextCell "E35" "  -6.47%  "
Set-TextCell "D36" "23.13"
Set-TextCell "E36" "  -4.56%  "
Set-TextCell "D37" "172.03"
Set-TextCell "E37" "  -0.60%  "
Set-TextCell "D38" "6.84"
Set-TextCell "E38" "  -9.33%  "
Set-TextCell "E39" "  -10.28%  "
Set-TextCell "D40" "1.47"
Set-TextCell "E40" "  -9.48%  "
Set-TextCell "E41" "  -8.61%  "
Set-TextCell "D42" "0.0781"
Set-TextCell "E42" "  -3.42%  "
Set-TextCell "D43" "0.813"
Set-TextCell "E43" "  -4.72%  "
Set-TextCell "B44" "OKB"
Set-TextCell "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D44" "42.56"
Set-TextCell "E44" "  -6.60%  "
Set-TextCell "B45" "FirstDigitalUSD"
Set-TextCell "C45" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D45" "0.999"
Set-TextCell "E45" "  +0.00%  "
Set-TextCell "D46" "25.02"
Set-TextCell "E46" "  +9.26%  "
Set-TextCell "D47" "4.36"
Set-TextCell "E47" "  -12.44%  "
Set-TextCell "D48" "1.63"
Set-TextCell "E48" "  -8.59%  "
Set-TextCell "D49" "1.16"
Set-TextCell "E49" "  +3.64%  "
Set-TextCell "D50" "6.69"
Set-TextCell "E50" "  -5.22%  "
Set-TextCell "B51" "dogwifhat"
Set-TextCell "C51" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D51" "2.07"
Set-TextCell "E51" "  -11.94%  "
